# Latest Specific Workspace and Add On Element
#
# Updates the Marketo workspace usage snapshot on Sheet1: swaps the
# workspace/model column values (Default/Automation -> the new
# ".Customer Support New Hire Workspace"/"ACT-SS" workspace+model pair),
# refreshes all the per-asset counts, updates the Tags/Integration/Account
# Name summary fields, and appends the three new add-on element rows at
# the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a plain number but must stay stored as TEXT
# (matching the rest of that column). Applying a Text number format before
# assigning the value keeps Excel from "helpfully" re-typing it as a number.
$forceTextCells = @("C9","D9","C10","C12","D12","C13","D13","C14","B20","B21")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# ---- Row 2 (Asset Data header row) ----
$ws.Range("C2").Value = ".Customer Support New Hire Workspace"
$ws.Range("D2").Value = "ACT-SS"
$ws.Range("E2").Value = ""

# ---- Row 3 (Emails) ----
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = ""

# ---- Row 4 (Forms) ----
$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = ""

# ---- Row 5 (Landing Pages) ----
$ws.Range("B5").Value = 71
$ws.Range("C5").Value = 64
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = ""

# ---- Row 6 (Images and Files) ----
$ws.Range("B6").Value = 582
$ws.Range("C6").Value = 581
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = ""

# ---- Row 7 (Snippets) ----
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = ""

# ---- Row 8 (Campaign Data header row) ----
$ws.Range("C8").Value = ".Customer Support New Hire Workspace"
$ws.Range("D8").Value = "ACT-SS"
$ws.Range("E8").Value = ""

# ---- Row 9 (All Triggered Campaigns) ----
$ws.Range("B9").Value = 40
$ws.Range("C9").Value = "37"
$ws.Range("D9").Value = "3"
$ws.Range("F9").Value = ""

# ---- Row 10 (Active Triggered Campaigns) ----
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = "0"
$ws.Range("F10").Value = ""

# ---- Row 11 (Batch Campaigns - Repeating Schedule) ----
$ws.Range("F11").Value = ""

# ---- Row 12 (All Batch Campaigns) ----
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = "82"
$ws.Range("D12").Value = "2"
$ws.Range("E12").Value = ""

# ---- Row 13 (All Campaigns) ----
$ws.Range("B13").Value = 128
$ws.Range("C13").Value = "123"
$ws.Range("D13").Value = "5"
$ws.Range("E13").Value = ""

# ---- Row 14 (Active Campaigns) ----
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = "0"
$ws.Range("E14").Value = ""

# ---- Row 15 (Database Data header row) ----
$ws.Range("C15").Value = ".Customer Support New Hire Workspace"
$ws.Range("D15").Value = "ACT-SS"

# ---- Row 16 (Segmentations) ----
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0

# ---- Row 17 (Leads) ----
$ws.Range("B17").Value = 198
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0

# ---- Row 18 (Program Data header row) ----
$ws.Range("C18").Value = ".Customer Support New Hire Workspace"
$ws.Range("D18").Value = "ACT-SS"

# ---- Row 19 (Models) ----
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1

# ---- Row 20 (Tags) ----
$ws.Range("B20").Value = "2688"

# ---- Row 21 (Integration) ----
$ws.Range("B21").Value = "74"

# ---- Row 23 (Account Name) ----
$ws.Range("B23").Value = "Pradyumna Sahoo"

# ---- Row 27 (Total WorkSpace) ----
$ws.Range("B27").Value = 2

# ---- New rows 30-32: additional add-on elements ----
$ws.Range("A30").Value = "Target Account Management"
$ws.Range("B30").Value = "True"
$ws.Range("A31").Value = "Predictive Content"
$ws.Range("B31").Value = "True"
$ws.Range("A32").Value = "Web Personalization"
$ws.Range("B32").Value = "True"
